# Insert a new data row right above the current row 18 (pushing the
# existing rows 18..50 down to 19..51) and populate it with the new
# record's values. This mirrors the author's edit: a new weekly price
# observation was added to the top of the historical list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 18:50 down by inserting a new row at 18.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new record.
$ws.Cells.Item(18, 1).Value = 7
$ws.Cells.Item(18, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(18, 3).Value = 'Ñuble'
$ws.Cells.Item(18, 4).Value = 44797
$ws.Cells.Item(18, 5).Value = 16
$ws.Cells.Item(18, 6).Value = 100112013
$ws.Cells.Item(18, 7).Value = 'Alcachofa'
$ws.Cells.Item(18, 8).Value = 'Argentina(o)'
$ws.Cells.Item(18, 9).Value = 'Primera'
$ws.Cells.Item(18, 10).Value = 120
$ws.Cells.Item(18, 11).Value = 15000
$ws.Cells.Item(18, 12).Value = 16000
$ws.Cells.Item(18, 13).Value = 15500
$ws.Cells.Item(18, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(18, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(18, 16).Value = 310
$ws.Cells.Item(18, 17).Value = 50
$ws.Cells.Item(18, 18).Value = 'Hortaliza'

# Match the date-formatted style used by the other rows in column D.
$ws.Cells.Item(18, 4).NumberFormat = $ws.Cells.Item(19, 4).NumberFormat
